$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5 (existing "Nia" record): Pajak_Terhutang, Tanggal_Jatuh_Tempo and
#     Pajak get re-entered as plain text instead of numbers/dates. ---

# E5: 70000 (number) -> "70000" (text)
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "70000"
$ws.Range("E5").Style = "Normal"

# F5: 46233 (date, formatted YYYY-MM-DD) -> "2026-07-30 00:00:00" (text, no special format)
$ws.Range("F5").Value = "2026-07-30 00:00:00"
$ws.Range("F5").Style = "Normal"

# G5: 70000 (number) -> "70000" (text)
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "70000"
$ws.Range("G5").Style = "Normal"

# --- Row 6 (new "Tiara" record) ---

# A6: NIK, kept as text so the long digit string isn't coerced to a number
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "2345678990112444"
$ws.Range("A6").Style = "Normal"

$ws.Range("B6").Value = "BG6701HI"
$ws.Range("C6").Value = "Tiara"
$ws.Range("D6").Value = "Palembang"

# E6 / G6: plain numbers
$ws.Range("E6").Value = 90000
$ws.Range("G6").Value = 90000

# F6: date value, formatted the same way as F5 used to be (YYYY-MM-DD)
$ws.Range("F6").Value = 46233
$ws.Range("F6").NumberFormat = "YYYY-MM-DD"
